$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Row -> (new Week label, new Seasonality Index value)
$updates = @(
    @{ Row = 2;  Week = "W05"; Seasonality = 0.93 },
    @{ Row = 3;  Week = "W06"; Seasonality = 1.14 },
    @{ Row = 4;  Week = "W07"; Seasonality = 1.07 },
    @{ Row = 5;  Week = "W08"; Seasonality = 1.06 },
    @{ Row = 6;  Week = "W09"; Seasonality = 0.99 },
    @{ Row = 7;  Week = "W10"; Seasonality = 1.05 },
    @{ Row = 8;  Week = "W11"; Seasonality = 1.17 },
    @{ Row = 9;  Week = "W12"; Seasonality = 1.03 },
    @{ Row = 10; Week = "W13"; Seasonality = 1.15 },
    @{ Row = 11; Week = "W14"; Seasonality = 1.17 },
    @{ Row = 12; Week = "W15"; Seasonality = 1.05 },
    @{ Row = 13; Week = "W16"; Seasonality = 0.91 },
    @{ Row = 14; Week = "W17"; Seasonality = 1.09 },
    @{ Row = 15; Week = "W18"; Seasonality = 1.03 },
    @{ Row = 16; Week = "W19"; Seasonality = 0.95 },
    @{ Row = 17; Week = "W20"; Seasonality = 1.03 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Week
    $ws.Cells.Item($u.Row, 16).Value = $u.Seasonality
}
